$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- First summary block (rows 15-25) ---

# New column header "charge [mC]" in K15, matching style of the other
# header cells in that row (same formatting as B15:J15).
$ws.Range("K15").Value = "charge [mC]"
$ws.Range("J15").Copy()
$ws.Range("K15").PasteSpecial(-4122)  # xlPasteFormats

# Row 16: "sleep" -> measured current changed from 50uA to 30uA
$ws.Range("E16").Value = 30
# New "charge [mC]" column for row 16
$ws.Range("K16").Formula = "=+D16*E16/1000"
# New measured-charge-per-wake figures off to the side (Q/R, like the
# Q3/R3 annotations used elsewhere on the sheet)
$ws.Range("Q16").Value = 87
$ws.Range("R16").Value = "mC per wake, measured"

# Row 17: "Active, Wifi on" -> duration changed from 300ms to 600ms
$ws.Range("D17").Value = 600
# New "charge [mC]" column for row 17
$ws.Range("K17").Formula = "=+D17*E17/1000"
# New calculated-charge-per-wake total next to it
$ws.Range("Q17").Formula = "=SUM(K16:K17)"
$ws.Range("R17").Value = "mC per wake, calculated"

# Row 19: "sleep" baseline current tweak 14uA -> 5uA
$ws.Range("F19").Value = 5

# --- Second summary block (rows 28-38), same underlying edits ---

# Row 29: same 50uA -> 30uA change
$ws.Range("E29").Value = 30

# Row 30: same 300ms -> 600ms change
$ws.Range("D30").Value = 600

# Row 32: same 14uA -> 5uA change
$ws.Range("F32").Value = 5

# --- View state: scroll up and move the selection the way the author
# left it after finishing these edits. ---
$aw = $excel.ActiveWindow
$aw.ScrollRow = 10
$aw.ScrollColumn = 1
$ws.Range("F33").Select() | Out-Null
